{"js": "// Insert the new \"BSD / I/O Kit\" reading-notes section right after the\n// standalone \"BSD\" heading paragraph that sits near the end of the document\n// (immediately before the trailing blank paragraph).\n//\n// The new content is 26 paragraphs (some of them blank separators). We\n// build them as one newline-joined string and hand it to\n// Paragraph.insertParagraph(), which the host splits on \"\\n\" into\n// sibling <w:p> elements inserted, in order, right after the anchor\n// paragraph.\nconst newParagraphTexts = [\n  \"BSD\u5927\u90e8\u5206\u8ddf\u539f\u7248\u4e00\u6837\uff0c\u4f46\u662f\u4e5f\u6709\u5c11\u90e8\u5206\u4e3a\u4e86\u80fd\u8ddfI/O kit\u548cMach\u517c\u5bb9\uff0c\u4ece\u800c\u9020\u6210\u4e86\u6539\u53d8\",\n  \"BSD\u5305\u62ec\u4ee5\u4e0b\u51e0\u4e2a\u90e8\u5206\",\n  \"BSD-style \u8fdb\u7a0b\u6a21\u578b\",\n  \"\u4fe1\u53f7\u91cf\",\n  \"\u7528\u6237id\u6743\u9650\",\n  \"POSIX API\",\n  \"\u5f02\u6b65IO\",\n  \"BSD-style\u7684system call\",\n  \"TCP/IP\u534f\u8bae\u6808\u548cBSDsocket\",\n  \"NKE\u7f51\u7edc\u6838\u5fc3\u6269\u5c55\uff08\u548cxnu\u517c\u5bb9\u7684\u90e8\u5206\uff09\",\n  \"VFS\",\n  \"System V\uff08\u6211\u4e2a\u4eba\u8ba4\u4e3a\u8fd9\u91cc\u9762\u5e94\u8be5\u6307\u7684\u662f\u8fdb\u7a0b\u95f4\u901a\u4fe1\u7684\u6807\u51c6\uff0c\u5305\u62ec\u6d88\u606f\u961f\u5217\uff0c\u5171\u4eab\u5185\u5b58\uff0c\u4fe1\u53f7\u91cf\uff09\u548cPOSIX\u8fdb\u7a0b\u95f4\u901a\u4fe1\u673a\u5236\",\n  \"\u5185\u6838\u90e8\u5206\u5bc6\u7801\u6846\u67b6\",\n  \"\u5df4\u62c9\u5df4\u62c9\",\n  \"\",\n  \"UBC\uff08unified buffer cache\uff09\u7edf\u4e00\u7f13\u5b58\u3002\",\n  \"\u8ba9\u6587\u4ef6\u53ef\u4ee5\u548c\u865a\u62df\u5185\u5b58\u4e00\u6837\u7f13\u5b58\u5728\u5185\u5b58\u4e2d\",\n  \"\",\n  \"I/O Kit\",\n  \"\u8fd9\u73a9\u610f\u662f\u7528C++\u5199\u7684\uff0c\u4f46\u662f\u9609\u5272\u4e86\u5f88\u591a\u7279\u6027\u4e0d\u8ba9\u7528\uff0c\u7136\u540e\u81ea\u5df1\u5b9e\u73b0\u4e86\u4e00\u4e2aRTTI\u7cfb\u7edf\",\n  \"\u5305\u62ec\u4e00\u4e2a\u5185\u6838\u5185\u7684C++\u5e93\uff0c\u548c\u4e00\u4e2a\u7528\u6237\u7a7a\u95f4\u7684\u6846\u67b6\",\n  \"I/O kit\u672c\u8eab\u662f\u6a21\u5757\u5316\u7684\u4ee5\u53ca\u5206\u5c42\u7684\uff0c\u5b83\u7ed9\u6355\u6349\uff0c\u8868\u793a\uff0c\u7ef4\u62a4\u6d89\u53ca\u5230I/O\u8054\u7cfb\u7684\u4e0d\u540c\u786c\u4ef6\u8f6f\u4ef6\u7ec4\u4ef6\u4e4b\u95f4\u7684\u5173\u7cfb\u63d0\u4f9b\u4e86\u4e00\u4e2a\u57fa\u7840\u8bbe\u65bd\u3002\",\n  \"\u4ed6\u5c01\u88c5\u4e86\u5e95\u5c42\u7684\u62bd\u8c61\u7ed9\u7cfb\u7edf\u7684\u5176\u4ed6\u90e8\u5206\u3002\",\n  \"\",\n  \"Libkern \u5e93\",\n  \"\u8fd9\u73a9\u610f\u5c31\u63d0\u4f9b\u4e86\u524d\u9762\u8bf4\u7684I/O kit\u7684runtime system\u3002\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the LAST paragraph whose entire text is exactly \"BSD\" \u2014 that's the\n// short standalone heading right before the final (empty) paragraph of the\n// document, which is where the new section gets appended.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"BSD\") {\n    anchor = paragraphs.items[i];\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find the trailing standalone \"BSD\" paragraph.');\n}\n\nanchor.insertParagraph(newParagraphTexts.join(\"\\n\"), Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert the new \"BSD / I/O Kit\" reading-notes section right after the\n# standalone \"BSD\" heading paragraph that sits near the end of the document\n# (immediately before the trailing blank paragraph).\n\n$newParagraphTexts = @(\n    'BSD\u5927\u90e8\u5206\u8ddf\u539f\u7248\u4e00\u6837\uff0c\u4f46\u662f\u4e5f\u6709\u5c11\u90e8\u5206\u4e3a\u4e86\u80fd\u8ddfI/O kit\u548cMach\u517c\u5bb9\uff0c\u4ece\u800c\u9020\u6210\u4e86\u6539\u53d8',\n    'BSD\u5305\u62ec\u4ee5\u4e0b\u51e0\u4e2a\u90e8\u5206',\n    'BSD-style \u8fdb\u7a0b\u6a21\u578b',\n    '\u4fe1\u53f7\u91cf',\n    '\u7528\u6237id\u6743\u9650',\n    'POSIX API',\n    '\u5f02\u6b65IO',\n    'BSD-style\u7684system call',\n    'TCP/IP\u534f\u8bae\u6808\u548cBSDsocket',\n    'NKE\u7f51\u7edc\u6838\u5fc3\u6269\u5c55\uff08\u548cxnu\u517c\u5bb9\u7684\u90e8\u5206\uff09',\n    'VFS',\n    'System V\uff08\u6211\u4e2a\u4eba\u8ba4\u4e3a\u8fd9\u91cc\u9762\u5e94\u8be5\u6307\u7684\u662f\u8fdb\u7a0b\u95f4\u901a\u4fe1\u7684\u6807\u51c6\uff0c\u5305\u62ec\u6d88\u606f\u961f\u5217\uff0c\u5171\u4eab\u5185\u5b58\uff0c\u4fe1\u53f7\u91cf\uff09\u548cPOSIX\u8fdb\u7a0b\u95f4\u901a\u4fe1\u673a\u5236',\n    '\u5185\u6838\u90e8\u5206\u5bc6\u7801\u6846\u67b6',\n    '\u5df4\u62c9\u5df4\u62c9',\n    '',\n    'UBC\uff08unified buffer cache\uff09\u7edf\u4e00\u7f13\u5b58\u3002',\n    '\u8ba9\u6587\u4ef6\u53ef\u4ee5\u548c\u865a\u62df\u5185\u5b58\u4e00\u6837\u7f13\u5b58\u5728\u5185\u5b58\u4e2d',\n    '',\n    'I/O Kit',\n    '\u8fd9\u73a9\u610f\u662f\u7528C++\u5199\u7684\uff0c\u4f46\u662f\u9609\u5272\u4e86\u5f88\u591a\u7279\u6027\u4e0d\u8ba9\u7528\uff0c\u7136\u540e\u81ea\u5df1\u5b9e\u73b0\u4e86\u4e00\u4e2aRTTI\u7cfb\u7edf',\n    '\u5305\u62ec\u4e00\u4e2a\u5185\u6838\u5185\u7684C++\u5e93\uff0c\u548c\u4e00\u4e2a\u7528\u6237\u7a7a\u95f4\u7684\u6846\u67b6',\n    'I/O kit\u672c\u8eab\u662f\u6a21\u5757\u5316\u7684\u4ee5\u53ca\u5206\u5c42\u7684\uff0c\u5b83\u7ed9\u6355\u6349\uff0c\u8868\u793a\uff0c\u7ef4\u62a4\u6d89\u53ca\u5230I/O\u8054\u7cfb\u7684\u4e0d\u540c\u786c\u4ef6\u8f6f\u4ef6\u7ec4\u4ef6\u4e4b\u95f4\u7684\u5173\u7cfb\u63d0\u4f9b\u4e86\u4e00\u4e2a\u57fa\u7840\u8bbe\u65bd\u3002',\n    '\u4ed6\u5c01\u88c5\u4e86\u5e95\u5c42\u7684\u62bd\u8c61\u7ed9\u7cfb\u7edf\u7684\u5176\u4ed6\u90e8\u5206\u3002',\n    '',\n    'Libkern \u5e93',\n    '\u8fd9\u73a9\u610f\u5c31\u63d0\u4f9b\u4e86\u524d\u9762\u8bf4\u7684I/O kit\u7684runtime system\u3002'\n)\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n# Find the LAST paragraph whose entire text is exactly \"BSD\" (ignoring the\n# trailing paragraph-mark character) -- the standalone heading right before\n# the final (empty) paragraph of the document.\n$targetIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $t = $paras.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"BSD\") {\n        $targetIndex = $i\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the trailing standalone 'BSD' paragraph.\"\n}\n\n$anchorRange = $paras.Item($targetIndex).Range\n$anchorRange.Collapse(0)\n$anchorRange.InsertParagraphAfter()\n$anchorRange.Collapse(0)\n\n$joined = [string]::Join([char]13, $newParagraphTexts)\n$anchorRange.InsertAfter($joined)\n\n$d.Save()\n"}
